$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (top header labels) ---
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"

# --- Rows 3-6: update anchor-score stats block (A-H) ---
$ws.Range("B3").Value = 0.9117647058823529
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 31
$ws.Range("H3").Value = 3

$ws.Range("B4").Value = 0.5924657534246576
$ws.Range("C4").Value = 173
$ws.Range("D4").Value = 173
$ws.Range("H4").Value = 119

$ws.Range("B5").Value = 0.1937984496124031
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 100
$ws.Range("H5").Value = 416

$ws.Range("B6").Value = 0.1746031746031746
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 33
$ws.Range("H6").Value = 156

# --- Row 7: clear the old "low" category row entirely (A-H) ---
$ws.Range("A7:H7").Clear()

# --- Rows 3-27: update J (word), K, L, M, Q for the ranked-word block ---
$ws.Range("J3").Value = "happy"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 26
$ws.Range("M3").Value = 26
$ws.Range("Q3").Value = 0
$ws.Range("J4").Value = "best"
$ws.Range("K4").Value = 0.9491525423728814
$ws.Range("L4").Value = 56
$ws.Range("M4").Value = 56
$ws.Range("Q4").Value = 3
$ws.Range("J5").Value = "interesting"
$ws.Range("K5").Value = 0.9090909090909091
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 30
$ws.Range("Q5").Value = 3
$ws.Range("J6").Value = "love"
$ws.Range("K6").Value = 0.8913043478260869
$ws.Range("L6").Value = 41
$ws.Range("M6").Value = 41
$ws.Range("Q6").Value = 5
$ws.Range("J7").Value = "great"
$ws.Range("K7").Value = 0.8482142857142857
$ws.Range("L7").Value = 95
$ws.Range("M7").Value = 95
$ws.Range("Q7").Value = 17
$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8170731707317073
$ws.Range("L8").Value = 67
$ws.Range("M8").Value = 67
$ws.Range("Q8").Value = 15
$ws.Range("J9").Value = "thank"
$ws.Range("K9").Value = 0.7890625
$ws.Range("L9").Value = 101
$ws.Range("M9").Value = 101
$ws.Range("Q9").Value = 27
$ws.Range("J10").Value = "positive"
$ws.Range("K10").Value = 0.7586206896551724
$ws.Range("L10").Value = 44
$ws.Range("M10").Value = 44
$ws.Range("Q10").Value = 14
$ws.Range("J11").Value = "free"
$ws.Range("K11").Value = 0.7333333333333333
$ws.Range("L11").Value = 88
$ws.Range("M11").Value = 88
$ws.Range("Q11").Value = 32
$ws.Range("J12").Value = "safe"
$ws.Range("K12").Value = 0.7183098591549296
$ws.Range("L12").Value = 102
$ws.Range("M12").Value = 102
$ws.Range("Q12").Value = 40
$ws.Range("J13").Value = "support"
$ws.Range("K13").Value = 0.7169811320754716
$ws.Range("L13").Value = 76
$ws.Range("M13").Value = 76
$ws.Range("Q13").Value = 30
$ws.Range("J14").Value = "good"
$ws.Range("K14").Value = 0.69375
$ws.Range("L14").Value = 111
$ws.Range("M14").Value = 111
$ws.Range("Q14").Value = 49
$ws.Range("J15").Value = "safety"
$ws.Range("K15").Value = 0.6862745098039216
$ws.Range("L15").Value = 35
$ws.Range("M15").Value = 35
$ws.Range("Q15").Value = 16
$ws.Range("J16").Value = "heroes"
$ws.Range("K16").Value = 0.6808510638297872
$ws.Range("L16").Value = 32
$ws.Range("M16").Value = 32
$ws.Range("Q16").Value = 15
$ws.Range("J17").Value = "relief"
$ws.Range("K17").Value = 0.66
$ws.Range("L17").Value = 33
$ws.Range("M17").Value = 33
$ws.Range("Q17").Value = 17
$ws.Range("J18").Value = "well"
$ws.Range("K18").Value = 0.6276595744680851
$ws.Range("L18").Value = 59
$ws.Range("M18").Value = 59
$ws.Range("Q18").Value = 35
$ws.Range("J19").Value = "better"
$ws.Range("K19").Value = 0.6031746031746031
$ws.Range("L19").Value = 38
$ws.Range("M19").Value = 38
$ws.Range("Q19").Value = 25
$ws.Range("J20").Value = "fresh"
$ws.Range("K20").Value = 0.5625
$ws.Range("L20").Value = 27
$ws.Range("M20").Value = 27
$ws.Range("Q20").Value = 21
$ws.Range("J21").Value = "hand"
$ws.Range("K21").Value = 0.5221932114882507
$ws.Range("L21").Value = 200
$ws.Range("M21").Value = 200
$ws.Range("Q21").Value = 183
$ws.Range("J22").Value = "help"
$ws.Range("K22").Value = 0.4915254237288136
$ws.Range("L22").Value = 145
$ws.Range("M22").Value = 145
$ws.Range("Q22").Value = 150
$ws.Range("J23").Value = "like"
$ws.Range("K23").Value = 0.4823529411764706
$ws.Range("L23").Value = 164
$ws.Range("M23").Value = 164
$ws.Range("Q23").Value = 176
$ws.Range("J24").Value = "care"
$ws.Range("K24").Value = 0.4606741573033708
$ws.Range("L24").Value = 41
$ws.Range("M24").Value = 41
$ws.Range("Q24").Value = 48
$ws.Range("J25").Value = "protect"
$ws.Range("K25").Value = 0.4520547945205479
$ws.Range("L25").Value = 33
$ws.Range("M25").Value = 33
$ws.Range("Q25").Value = 40
$ws.Range("J26").Value = "increase"
$ws.Range("K26").Value = 0.4102564102564102
$ws.Range("L26").Value = 32
$ws.Range("M26").Value = 32
$ws.Range("Q26").Value = 46
$ws.Range("J27").Value = "please"
$ws.Range("K27").Value = 0.3389121338912134
$ws.Range("L27").Value = 81
$ws.Range("M27").Value = 81
$ws.Range("Q27").Value = 158

# --- Row 28: now also N/O/P change (became a below-threshold "both" row) ---
$ws.Range("J28").Value = "store"
$ws.Range("K28").Value = 0.04031354983202688
$ws.Range("L28").Value = 36
$ws.Range("M28").Value = 37
$ws.Range("N28").Value = 0.97
$ws.Range("O28").Value = 0.03000000000000003
$ws.Range("P28").Value = $true
$ws.Range("Q28").Value = 857

# --- Rows 29-30: brand new rows, copy formatting from row 28's J cell ---
$ws.Range("J28").Copy()
$ws.Range("J29:J30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("J29").Value = "grocery"
$ws.Range("K29").Value = 0.03107658157602664
$ws.Range("L29").Value = 28
$ws.Range("M29").Value = 28
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 873

$ws.Range("J30").Value = "co"
$ws.Range("K30").Value = 0.00933977455716586
$ws.Range("L30").Value = 29
$ws.Range("M30").Value = 31
$ws.Range("N30").Value = 0.94
$ws.Range("O30").Value = 0.06000000000000005
$ws.Range("P30").Value = $true
$ws.Range("Q30").Value = 3076

